# Commit: "rD and rs1 should be all possible GPRs."
#
# On the "MikeSuggestedFormat" worksheet, prepend a note about rD/rs1
# needing to cover all GPRs to the "Verification Goal" cells (column E)
# of the four extract-instruction rows (2-5), grow those rows so the
# extra line still fits, and leave the selection on A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MikeSuggestedFormat")
$ws.Activate()

$note = "rD and rs1 should be all possible GPRs.`n"

foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $note + $cell.Value2
    $ws.Rows.Item($r).RowHeight = 69
}

$ws.Range("A2").Select()
